# Update the answers in the "two-digit divided by one-digit" practice
# table. The document contains a single table; data values live in
# Word table rows 1, 5, 9, 13 and 17 (the other rows are blank spacer
# rows), five columns each.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Map of (row, col) -> new text, row/col are 1-based Word table indices.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "37÷7=5, 2" },
    @{ Row = 1;  Col = 2; Text = "75÷2=37, 1" },
    @{ Row = 1;  Col = 3; Text = "10÷3=3, 1" },
    @{ Row = 1;  Col = 4; Text = "87÷9=9, 6" },
    @{ Row = 1;  Col = 5; Text = "92÷7=13, 1" },

    @{ Row = 5;  Col = 1; Text = "44÷2=22, 0" },
    @{ Row = 5;  Col = 2; Text = "68÷3=22, 2" },
    @{ Row = 5;  Col = 3; Text = "19÷5=3, 4" },
    @{ Row = 5;  Col = 4; Text = "39÷5=7, 4" },
    @{ Row = 5;  Col = 5; Text = "33÷9=3, 6" },

    @{ Row = 9;  Col = 1; Text = "16÷7=2, 2" },
    @{ Row = 9;  Col = 2; Text = "75÷2=37, 1" },
    @{ Row = 9;  Col = 3; Text = "18÷3=6, 0" },
    @{ Row = 9;  Col = 4; Text = "43÷5=8, 3" },
    @{ Row = 9;  Col = 5; Text = "29÷2=14, 1" },

    @{ Row = 13; Col = 1; Text = "44÷9=4, 8" },
    @{ Row = 13; Col = 2; Text = "57÷7=8, 1" },
    @{ Row = 13; Col = 3; Text = "62÷4=15, 2" },
    @{ Row = 13; Col = 4; Text = "68÷5=13, 3" },
    @{ Row = 13; Col = 5; Text = "64÷3=21, 1" },

    @{ Row = 17; Col = 1; Text = "35÷4=8, 3" },
    @{ Row = 17; Col = 2; Text = "67÷5=13, 2" },
    @{ Row = 17; Col = 3; Text = "84÷6=14, 0" },
    @{ Row = 17; Col = 4; Text = "38÷3=12, 2" },
    @{ Row = 17; Col = 5; Text = "55÷7=7, 6" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text of the cell is replaced.
    $rng.MoveEnd(12, -1) | Out-Null
    $rng.Text = $u.Text
}
